$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 45455136
$ws.Range("I80").Value = 362.4
$ws.Range("J80").Value = 83334120
$ws.Range("K80").Value = 1087.2
$ws.Range("L80").Value = 250002360
$ws.Range("M80").Value = -89.19999999999982
$ws.Range("N80").Value = -250004356
$ws.Range("H83").Value = 45455136
$ws.Range("I83").Value = 362.4
$ws.Range("J83").Value = 83334120
$ws.Range("K83").Value = 3261.6
$ws.Range("L83").Value = 750007080
$ws.Range("M83").Value = 1730.4
$ws.Range("N83").Value = -750017064
$ws.Range("H112").Value = 1715.6154
$ws.Range("J112").Value = 1941.6364
$ws.Range("L112").Value = 5824.9092
$ws.Range("N112").Value = -8040.9092
$ws.Range("H125").Value = 1755.1818
$ws.Range("J125").Value = 3891
$ws.Range("L125").Value = 35019
$ws.Range("M125").Value = -2352.428699999999
$ws.Range("N125").Value = -39939
$ws.Range("H127").Value = 1325.6471
$ws.Range("I127").Value = 306.33334
$ws.Range("J127").Value = 1544.0714
$ws.Range("K127").Value = 919.0000200000001
$ws.Range("L127").Value = 4632.2142
$ws.Range("M127").Value = 4040.99998
$ws.Range("N127").Value = -14552.2142
$ws.Range("H129").Value = 783.4706
$ws.Range("I129").Value = 403.91666
$ws.Range("J129").Value = 990.5
$ws.Range("K129").Value = 1211.74998
$ws.Range("L129").Value = 2971.5
$ws.Range("M129").Value = 3788.25002
$ws.Range("N129").Value = -12971.5
$ws.Range("H138").Value = 2384999
$ws.Range("I138").Value = 2004.6428
$ws.Range("J138").Value = 3973662
$ws.Range("K138").Value = 6013.928400000001
$ws.Range("L138").Value = 11920986
$ws.Range("M138").Value = -873.9284000000007
$ws.Range("N138").Value = -11931266
$ws.Range("H141").Value = 2846
$ws.Range("I141").Value = 2015.8334
$ws.Range("J141").Value = 6166.6665
$ws.Range("K141").Value = 6047.5002
$ws.Range("L141").Value = 18499.9995
$ws.Range("M141").Value = -867.5002000000004
$ws.Range("N141").Value = -28859.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4585475.5
$ws.Range("I32").Value = 5669361.5
$ws.Range("K32").Value = 5669361.5
$ws.Range("M32").Value = -5669074.5
$ws.Range("H122").Value = 55560560
$ws.Range("I122").Value = 10012
$ws.Range("J122").Value = 111111110
$ws.Range("K122").Value = 30036
$ws.Range("L122").Value = 333333330
$ws.Range("M122").Value = -27586
$ws.Range("N122").Value = -333338230
$ws.Range("H132").Value = 288642.56
$ws.Range("I132").Value = 252625
$ws.Range("J132").Value = 336666
$ws.Range("K132").Value = 757875
$ws.Range("L132").Value = 1009998
$ws.Range("M132").Value = -755345
$ws.Range("N132").Value = -1015058

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1170.4286
$ws.Range("I94").Value = 419.5
$ws.Range("J94").Value = 1470.8
$ws.Range("K94").Value = 419.5
$ws.Range("L94").Value = 1470.8
$ws.Range("M94").Value = 31.5
$ws.Range("N94").Value = -2372.8
$ws.Range("H107").Value = 2766.8333
$ws.Range("I107").Value = 2500.5
$ws.Range("K107").Value = 2500.5
$ws.Range("M107").Value = -580.5
$ws.Range("H138").Value = 30000
$ws.Range("J138").Value = 30000
$ws.Range("L138").Value = 30000
$ws.Range("N138").Value = -40280

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 47782530
$ws.Range("I6").Value = 63706704
$ws.Range("J6").Value = 9999
$ws.Range("K6").Value = 63706704
$ws.Range("L6").Value = 9999
$ws.Range("M6").Value = -63706591
$ws.Range("N6").Value = -10225
$ws.Range("H7").Value = 180.33333
$ws.Range("I7").Value = 40
$ws.Range("J7").Value = 250.5
$ws.Range("K7").Value = 40
$ws.Range("L7").Value = 250.5
$ws.Range("M7").Value = 73
$ws.Range("N7").Value = -476.5
$ws.Range("H31").Value = 8962.040999999999
$ws.Range("I31").Value = 22885.28
$ws.Range("J31").Value = 1710.3541
$ws.Range("K31").Value = 22885.28
$ws.Range("L31").Value = 1710.3541
$ws.Range("M31").Value = -22590.28
$ws.Range("N31").Value = -2300.3541
$ws.Range("H34").Value = 8962.040999999999
$ws.Range("I34").Value = 22885.28
$ws.Range("J34").Value = 1710.3541
$ws.Range("K34").Value = 22885.28
$ws.Range("L34").Value = 1710.3541
$ws.Range("M34").Value = -22683.28
$ws.Range("N34").Value = -2114.3541
$ws.Range("H58").Value = 17896938
$ws.Range("I58").Value = 23723318
$ws.Range("J58").Value = 1625.5714
$ws.Range("K58").Value = 23723318
$ws.Range("L58").Value = 1625.5714
$ws.Range("M58").Value = -23723115
$ws.Range("N58").Value = -2031.5714
$ws.Range("H87").Value = 36100
$ws.Range("J87").Value = 36100
$ws.Range("L87").Value = 36100
$ws.Range("N87").Value = -38472
$ws.Range("H90").Value = 36100
$ws.Range("J90").Value = 36100
$ws.Range("L90").Value = 108300
$ws.Range("N90").Value = -120156
$ws.Range("H134").Value = 24975.088
$ws.Range("I134").Value = 1410.9678
$ws.Range("J134").Value = 73674.266
$ws.Range("K134").Value = 4232.903399999999
$ws.Range("L134").Value = 221022.798
$ws.Range("M134").Value = -1697.903399999999
$ws.Range("N134").Value = -226092.798
$ws.Range("H136").Value = 17896938
$ws.Range("I136").Value = 23723318
$ws.Range("J136").Value = 1625.5714
$ws.Range("K136").Value = 71169954
$ws.Range("L136").Value = 4876.7142
$ws.Range("M136").Value = -71167404
$ws.Range("N136").Value = -9976.7142

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 10280.667
$ws.Range("I80").Value = 36083
$ws.Range("J80").Value = 3830.0833
$ws.Range("K80").Value = 36083
$ws.Range("L80").Value = 3830.0833
$ws.Range("M80").Value = -35085
$ws.Range("N80").Value = -5826.0833
$ws.Range("H83").Value = 10280.667
$ws.Range("I83").Value = 36083
$ws.Range("J83").Value = 3830.0833
$ws.Range("K83").Value = 180415
$ws.Range("L83").Value = 19150.4165
$ws.Range("M83").Value = -175423
$ws.Range("N83").Value = -29134.4165
$ws.Range("H122").Value = 3334.8333
$ws.Range("I122").Value = 3001.8
$ws.Range("K122").Value = 9005.400000000001
$ws.Range("M122").Value = -6555.400000000001
$ws.Range("H132").Value = 121695.65
$ws.Range("I132").Value = 115014.89
$ws.Range("J132").Value = 129211.5
$ws.Range("K132").Value = 345044.67
$ws.Range("L132").Value = 387634.5
$ws.Range("M132").Value = -342514.67
$ws.Range("N132").Value = -392694.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2800.75
$ws.Range("I61").Value = 2950.5
$ws.Range("J61").Value = 2501.25
$ws.Range("K61").Value = 2950.5
$ws.Range("L61").Value = 2501.25
$ws.Range("M61").Value = -2748.5
$ws.Range("N61").Value = -2905.25
$ws.Range("H93").Value = 2219.2856
$ws.Range("I93").Value = 2190.1667
$ws.Range("J93").Value = 2394
$ws.Range("K93").Value = 2190.1667
$ws.Range("L93").Value = 2394
$ws.Range("M93").Value = -942.1667000000002
$ws.Range("N93").Value = -4890
$ws.Range("H113").Value = 2800.75
$ws.Range("I113").Value = 2950.5
$ws.Range("J113").Value = 2501.25
$ws.Range("K113").Value = 2950.5
$ws.Range("L113").Value = 2501.25
$ws.Range("M113").Value = -780.5
$ws.Range("N113").Value = -6841.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 2100.25
$ws.Range("I5").Value = 1200.5
$ws.Range("J5").Value = 3000
$ws.Range("K5").Value = 1200.5
$ws.Range("L5").Value = 3000
$ws.Range("M5").Value = -1088.5
$ws.Range("N5").Value = -3224
$ws.Range("H96").Value = 10000
$ws.Range("I96").Value = 7666.6665
$ws.Range("J96").Value = 17000
$ws.Range("K96").Value = 7666.6665
$ws.Range("L96").Value = 17000
$ws.Range("M96").Value = -6293.6665
$ws.Range("N96").Value = -19746
$ws.Range("H100").Value = 94272.91
$ws.Range("I100").Value = 102900.4
$ws.Range("J100").Value = 87083.336
$ws.Range("K100").Value = 205800.8
$ws.Range("L100").Value = 174166.672
$ws.Range("M100").Value = -205259.8
$ws.Range("N100").Value = -175248.672
$ws.Range("H107").Value = 518.1111
$ws.Range("I107").Value = 457.875
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 1373.625
$ws.Range("L107").Value = 3000
$ws.Range("M107").Value = 546.375
$ws.Range("N107").Value = -6840
$ws.Range("H132").Value = 50554.195
$ws.Range("I132").Value = 41655.04
$ws.Range("J132").Value = 64459.125
$ws.Range("K132").Value = 124965.12
$ws.Range("L132").Value = 193377.375
$ws.Range("M132").Value = -122435.12
$ws.Range("N132").Value = -198437.375
$ws.Range("H136").Value = 55938.703
$ws.Range("I136").Value = 35397.83
$ws.Range("J136").Value = 130399.375
$ws.Range("K136").Value = 106193.49
$ws.Range("L136").Value = 391198.125
$ws.Range("M136").Value = -103643.49
$ws.Range("N136").Value = -396298.125
